# Fixed update to excel issue
# - Rename "Requested quantity" headers to metric-specific names
# - Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# Rename the "Requested quantity" header on each existing sheet
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new worksheet after the last existing sheet and name it "PO Forecast"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Range("A2").Value = 45151.99999999999
$wsForecast.Range("B2").Value = 20
$wsForecast.Range("C2").Value = 19.99984888231002
$wsForecast.Range("D2").Value = 19.9998488828808

$wsForecast.Range("A3").Value = 45158.99999999999
$wsForecast.Range("B3").Value = 25
$wsForecast.Range("C3").Value = 24.99984888490909
$wsForecast.Range("D3").Value = 24.9998488854833

$wsForecast.Range("A4").Value = 45165.99999999999
$wsForecast.Range("B4").Value = 30
$wsForecast.Range("C4").Value = 29.999848683796
$wsForecast.Range("D4").Value = 29.99984910002066

$wsForecast.Range("A5").Value = 45172.99999999999
$wsForecast.Range("B5").Value = 35
$wsForecast.Range("C5").Value = 34.9998482587991
$wsForecast.Range("D5").Value = 34.99984951671681

$wsForecast.Range("A6").Value = 45179.99999999999
$wsForecast.Range("B6").Value = 40
$wsForecast.Range("C6").Value = 39.99984770967628
$wsForecast.Range("D6").Value = 39.99985008999176

$wsForecast.Range("A7").Value = 45186.99999999999
$wsForecast.Range("B7").Value = 45
$wsForecast.Range("C7").Value = 44.99984694312253
$wsForecast.Range("D7").Value = 44.99985084178653

$wsForecast.Range("A8").Value = 45193.99999999999
$wsForecast.Range("B8").Value = 50
$wsForecast.Range("C8").Value = 49.99984609174802
$wsForecast.Range("D8").Value = 49.99985172375896

$wsForecast.Range("A9").Value = 45200.99999999999
$wsForecast.Range("B9").Value = 55
$wsForecast.Range("C9").Value = 54.99984523961517
$wsForecast.Range("D9").Value = 54.99985263244913

$wsForecast.Range("A10").Value = 45207.99999999999
$wsForecast.Range("B10").Value = 60
$wsForecast.Range("C10").Value = 59.99984431749051
$wsForecast.Range("D10").Value = 59.99985361223604

$wsForecast.Range("A11").Value = 45214.99999999999
$wsForecast.Range("B11").Value = 65
$wsForecast.Range("C11").Value = 64.99984320841598
$wsForecast.Range("D11").Value = 64.99985471828356

# Match formatting of the source sheets: bold/bordered/centered header style,
# and the date-time number format used for the "ds" column.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Select()
$excel.CutCopyMode = 0
